# Auto-generated edit script: updates crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.707.21'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '3.332.54'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.21'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.83'
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('D9').Value = '3.328.69'
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('E10').Value = '  +1.83%  '
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.35'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '703.91'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '3.881.33'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.43'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '67.726.88'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '3.334.73'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.893'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.42'
$ws.Range('E23').Value = '  +4.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.97'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.29'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('E26').Value = '  -1.70%  '
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.48'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.25'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.54'
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.10'
$ws.Range('E31').Value = '  +4.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '570.84'
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.29'
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('D37').Value = '3.705.36'
$ws.Range('E37').Value = '  -4.25%  '
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('E39').Value = '  +5.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.131'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('D43').Value = '0.0₃0675'
$ws.Range('E43').Value = '  -0.96%  '
$ws.Range('E44').Value = '  +1.44%  '
$ws.Range('E45').Value = '  -3.25%  '
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.67'
$ws.Range('E47').Value = '  +5.88%  '
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.32'
$ws.Range('E50').Value = '  -5.54%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.70'
$ws.Range('E51').Value = '  +15.08%  '
